$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10, shifting existing rows 10-99 down to 11-100.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new data record.
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44817
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112031
$ws.Range("G10").Value = "Poroto verde"
$ws.Range("H10").Value = "Magnum"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 35
$ws.Range("K10").Value = 37000
$ws.Range("L10").Value = 37000
$ws.Range("M10").Value = 37000
$ws.Range("N10").Value = "`$/malla 25 kilos"
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 1480
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"

# Ensure the date column keeps the same date/time number format used by the
# other rows in column D.
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
